$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H7").Value = 5583.3335
$ws.Range("J7").Value = 8500
$ws.Range("L7").Value = 8500
$ws.Range("N7").Value = -8724

$ws.Range("H14").Value = 5583.3335
$ws.Range("J14").Value = 8500
$ws.Range("L14").Value = 8500
$ws.Range("N14").Value = -8882

$ws.Range("H17").Value = 2043.5
$ws.Range("J17").Value = 2043.5
$ws.Range("L17").Value = 6130.5
$ws.Range("N17").Value = -6466.5

$ws.Range("H40").Value = 4753.5557
$ws.Range("I40").Value = 4449.5
$ws.Range("J40").Value = 4996.8
$ws.Range("K40").Value = 4449.5
$ws.Range("L40").Value = 4996.8
$ws.Range("M40").Value = -4274.5
$ws.Range("N40").Value = -5346.8

$ws.Range("H86").Value = 8332.666999999999
$ws.Range("I86").Value = 8332.666999999999
$ws.Range("K86").Value = 8332.666999999999
$ws.Range("M86").Value = -7209.666999999999

$ws.Range("H89").Value = 8332.666999999999
$ws.Range("I89").Value = 8332.666999999999
$ws.Range("K89").Value = 41663.335
$ws.Range("M89").Value = -36047.335

$ws.Range("H132").Value = 6199.5
$ws.Range("J132").Value = 9895
$ws.Range("L132").Value = 29685
$ws.Range("N132").Value = -34745

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 6907.3896
$ws.Range("I32").Value = 5326.367
$ws.Range("J32").Value = 14654.4
$ws.Range("K32").Value = 5326.367
$ws.Range("L32").Value = 14654.4
$ws.Range("M32").Value = -5039.367
$ws.Range("N32").Value = -15228.4

$ws.Range("H45").Value = 3569.76
$ws.Range("I45").Value = 963.5714
$ws.Range("K45").Value = 963.5714
$ws.Range("M45").Value = -586.5714

$ws.Range("H74").Value = 1241.1765
$ws.Range("I74").Value = 828.36365
$ws.Range("K74").Value = 828.36365
$ws.Range("M74").Value = 45.63634999999999

$ws.Range("H77").Value = 1241.1765
$ws.Range("I77").Value = 828.36365
$ws.Range("K77").Value = 4141.81825
$ws.Range("M77").Value = 226.1817499999997

$ws.Range("H122").Value = 4116.909
$ws.Range("I122").Value = 4652.606
$ws.Range("J122").Value = 2509.818
$ws.Range("K122").Value = 13957.818
$ws.Range("L122").Value = 7529.454000000001
$ws.Range("M122").Value = -11507.818
$ws.Range("N122").Value = -12429.454

$ws.Range("H132").Value = 4285.844
$ws.Range("I132").Value = 4401.225
$ws.Range("K132").Value = 13203.675
$ws.Range("M132").Value = -10673.675

$ws.Range("H134").Value = 79950
$ws.Range("J134").Value = 79950
$ws.Range("L134").Value = 79950
$ws.Range("N134").Value = -90090

$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 719.15
$ws.Range("I94").Value = 386.75
$ws.Range("K94").Value = 386.75
$ws.Range("M94").Value = 64.25

$ws.Range("H134").Value = 7137.659
$ws.Range("I134").Value = 5581.8667
$ws.Range("K134").Value = 16745.6001
$ws.Range("M134").Value = -14210.6001

$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 226.75
$ws.Range("I7").Value = 185.84616
$ws.Range("K7").Value = 185.84616
$ws.Range("M7").Value = -72.84616

$ws.Range("H99").Value = 2222.2
$ws.Range("I99").Value = 799.5
$ws.Range("J99").Value = 3170.6667
$ws.Range("K99").Value = 799.5
$ws.Range("L99").Value = 3170.6667
$ws.Range("M99").Value = 698.5
$ws.Range("N99").Value = -6166.6667

$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 400
$ws.Range("M107").Value = 1520

$ws.Range("H126").Value = 2222.2
$ws.Range("I126").Value = 799.5
$ws.Range("J126").Value = 3170.6667
$ws.Range("K126").Value = 2398.5
$ws.Range("L126").Value = 9512.000100000001
$ws.Range("M126").Value = 71.5
$ws.Range("N126").Value = -14452.0001

$ws.Range("H141").Value = 157084.1
$ws.Range("J141").Value = 251820.83
$ws.Range("L141").Value = 251820.83
$ws.Range("N141").Value = -262180.83

$ws = $wb.Worksheets.Item(5)
$ws.Range("H48").Value = 2995
$ws.Range("I48").Value = 2995
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 8985
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -8735
$ws.Range("N48").ClearContents()

$ws.Range("H138").Value = 12845.208
$ws.Range("I138").Value = 3057
$ws.Range("K138").Value = 9171
$ws.Range("M138").Value = -4031

$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 11314.091
$ws.Range("I70").Value = 13776.833
$ws.Range("K70").Value = 13776.833
$ws.Range("M70").Value = -13506.833

$ws.Range("H73").Value = 11314.091
$ws.Range("I73").Value = 13776.833
$ws.Range("K73").Value = 13776.833
$ws.Range("M73").Value = -12840.833

$ws.Range("H80").Value = 3937.25
$ws.Range("I80").Value = 3599.8
$ws.Range("J80").Value = 4499.6665
$ws.Range("K80").Value = 3599.8
$ws.Range("L80").Value = 4499.6665
$ws.Range("M80").Value = -2601.8
$ws.Range("N80").Value = -6495.6665

$ws.Range("H83").Value = 3937.25
$ws.Range("I83").Value = 3599.8
$ws.Range("J83").Value = 4499.6665
$ws.Range("K83").Value = 17999
$ws.Range("L83").Value = 22498.3325
$ws.Range("M83").Value = -13007
$ws.Range("N83").Value = -32482.3325

$ws.Range("H122").Value = 34701.824
$ws.Range("I122").Value = 52047.19
$ws.Range("J122").Value = 6682.385
$ws.Range("K122").Value = 156141.57
$ws.Range("L122").Value = 20047.155
$ws.Range("M122").Value = -153691.57
$ws.Range("N122").Value = -24947.155

$ws.Range("H126").Value = 9293
$ws.Range("I126").Value = 9799.799999999999
$ws.Range("J126").Value = 8870.666999999999
$ws.Range("K126").Value = 29399.4
$ws.Range("L126").Value = 26612.001
$ws.Range("M126").Value = -26929.4
$ws.Range("N126").Value = -31552.001

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 9864.166999999999
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H16").Value = 1218.9333
$ws.Range("I16").Value = 1170.2858
$ws.Range("K16").Value = 1170.2858
$ws.Range("M16").Value = -1000.2858

$ws.Range("H40").Value = 4299
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H82").Value = 4022
$ws.Range("I82").Value = 1366.4
$ws.Range("J82").Value = 5349.8
$ws.Range("K82").Value = 1366.4
$ws.Range("L82").Value = 5349.8
$ws.Range("M82").Value = -1005.4
$ws.Range("N82").Value = -6071.8

$ws.Range("H85").Value = 4022
$ws.Range("I85").Value = 1366.4
$ws.Range("J85").Value = 5349.8
$ws.Range("K85").Value = 1366.4
$ws.Range("L85").Value = 5349.8
$ws.Range("M85").Value = -118.4000000000001
$ws.Range("N85").Value = -7845.8

$ws.Range("H126").Value = 9864.166999999999
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 3639.5588
$ws.Range("I136").Value = 3685.524
$ws.Range("K136").Value = 11056.572
$ws.Range("M136").Value = -8506.572

$ws = $wb.Worksheets.Item(8)
$ws.Range("H45").Value = 38631.777
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 38631.777
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 38631.777
$ws.Range("N45").Value = -39613.777
$ws.Range("M45").ClearContents()

$ws.Range("H132").Value = 3003.6538
$ws.Range("I132").Value = 2670.625
$ws.Range("K132").Value = 8011.875
$ws.Range("M132").Value = -5481.875

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
